$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new row of data (row 21): date, hours, and activity description.
# Copy the date style used by the rows above (s="4", numFmtId 14) instead of
# minting a new number format, and keep the value a pure date serial (no time).
$ws.Range("A20").Copy()
$ws.Range("A21").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A21").Value = 44246
$ws.Range("B21").Value = 4
$ws.Range("D21").Value = "Development WebUI"

# Recalculate so the running-sum formulas in column C pick up the new value.
$excel.Calculate()

# Update the selection to match what was last selected in the saved file.
$ws.Range("E26").Select()
